$wb = $excel.ActiveWorkbook

# Update "展览" sheet (row 2 -> F2: 169 -> 172, row 5 -> F5: 42 -> 43)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 172
$ws1.Range("F5").Value = 43

# Update "全部类型" sheet (same data mirrored, same changes)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 172
$ws4.Range("F5").Value = 43
